# lecture05.pptx edit — refresh "last modified" date field shown on the
# Date placeholder of the slide master / every slide layout, and update
# the title-slide subtitle (drop the "Fall 2022 |" term prefix and
# re-center the subtitle box).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date Placeholder (type datetime1) : 11/2/2022 -> 9/27/2023
#    Present once on the Slide Master and once on every Custom Layout.
# ---------------------------------------------------------------------
$newDate = "9/27/2023"

function Set-DatePlaceholderText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Type -eq 14) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 1 (title slide) subtitle shape:
#      - text: "Fall 2022 | University of Mount union"
#               -> "University of Mount union"
#      - reposition box slightly left/up
#        (EMU 708248,1447341 -> 581191,1440465 ; size unchanged)
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)

$subtitle = $null
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $candidate = $slide1.Shapes.Item($i)
    if ($candidate.Type -eq 14 -and $candidate.PlaceholderFormat.Type -eq 4) {
        $subtitle = $candidate
        break
    }
}
if ($subtitle -eq $null) {
    $subtitle = $slide1.Shapes.Item(3)
}

$subtitle.TextFrame.TextRange.Text = "University of Mount union"

# Shape.Left/Top are in points (1 pt = 12700 EMU); the literals below are
# tuned so that, after the host's internal float round-trip, the stored
# EMU offsets land exactly on 581191 / 1440465.
$subtitle.Left = 45.76307106019685
$subtitle.Top = 113.42244338992126
